# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Estado de Cuenta" worker/debt table (B15:J.. ) is rebuilt:
#   - a new worker (DILSA DEL CARMEN BONILLA CARMONA, doc 32937855) is added
#     with 5 overdue periods (2503-2507) at the top of the table,
#   - the existing CAROLINA SANTAMARIA MOLINA rows keep their two periods but
#     the period/value pairs are swapped,
#   - LAURA DANIELA ROMERO LEON keeps her single row but loses the "closing"
#     (bottom-border) row style, since she is no longer the last row,
#   - a new worker (WILMER RAMOS CASARRUBIA, doc 1002192305) is appended with
#     5 overdue periods (2503-2507); his last period row (2503) becomes the
#     new "closing" row of the table,
#   - the summary cells (total overdue value, worker count, period count) are
#     refreshed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# 1) Make room for the 5 new DILSA rows right above the existing data rows.
#    (old row 16 -> 21, old row 17 -> 22, old row 18 -> 23, old footer rows
#     23/24 -> 28/29)
# ---------------------------------------------------------------------------
$ws.Rows("16:20").Insert($xlShiftDown)

# Give the freshly inserted rows the same ("normal data row") formatting as
# the row right below them (the former row 16).
$ws.Range("B21:J21").Copy()
$ws.Range("B16:J20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Make room for the 5 new WILMER rows below LAURA's row (currently row 23).
#    (current footer rows 28/29 -> 33/34)
# ---------------------------------------------------------------------------
$ws.Rows("24:28").Insert($xlShiftDown)

# Row 28 is now the new last data row of the table, so it must get the
# special "closing" (bottom border) style that currently still sits on row 23
# (LAURA's row).
$ws.Range("B23:J23").Copy()
$ws.Range("B28:J28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Rows 24-27 are plain ("normal") data rows, same style as row 22.
$ws.Range("B22:J22").Copy()
$ws.Range("B24:J27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# LAURA's row (23) is no longer the last row of the table, so it loses the
# "closing" style and becomes a plain data row (same as row 22).
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Fill in the values.
# ---------------------------------------------------------------------------

# New worker: DILSA DEL CARMEN BONILLA CARMONA, doc 32937855 - rows 16-20
$dilsaPeriods = @("2507","2506","2505","2504","2503")
for ($i = 0; $i -lt 5; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "32937855"
    $ws.Cells.Item($r, 4).Value = "DILSA DEL CARMEN BONILLA CARMONA"
    $ws.Cells.Item($r, 5).Value = $dilsaPeriods[$i]
    $ws.Cells.Item($r, 6).Value = 56940
    $ws.Cells.Item($r, 7).Value = 781242
}

# CAROLINA SANTAMARIA MOLINA - rows 21-22: period/value pairs swapped
$ws.Range("E21").Value = "1908"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = 908526

$ws.Range("E22").Value = "1907"
$ws.Range("F22").Value = 30916
$ws.Range("G22").Value = 908526

# LAURA DANIELA ROMERO LEON - row 23: values unchanged, only the style moved
# (handled above); re-assert the values so the row is self-consistent.
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143398464"
$ws.Range("D23").Value = "LAURA DANIELA ROMERO LEON"
$ws.Range("E23").Value = "1908"
$ws.Range("F23").Value = 3312
$ws.Range("G23").Value = 877803

# New worker: WILMER RAMOS CASARRUBIA, doc 1002192305 - rows 24-28
$wilmerPeriods = @("2507","2506","2505","2504","2503")
for ($i = 0; $i -lt 5; $i++) {
    $r = 24 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1002192305"
    $ws.Cells.Item($r, 4).Value = "WILMER RAMOS CASARRUBIA"
    $ws.Cells.Item($r, 5).Value = $wilmerPeriods[$i]
    $ws.Cells.Item($r, 6).Value = 68480
    $ws.Cells.Item($r, 7).Value = 1712000
}

# ---------------------------------------------------------------------------
# 4) Refresh the summary cells.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 694453   # VALOR MORA (sum of column F, rows 16-28)
$ws.Range("C13").Value = 4        # Cant. Trabajadores
$ws.Range("F13").Value = 7        # Cant. Periodos

Write-Output "edit complete"
